$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 265; everything from 265..319 shifts down to 266..320.
$ws.Rows(265).Insert()

# Populate the newly inserted row 265 with the new weekly data point.
$ws.Cells.Item(265, 1).Value = 11
$ws.Cells.Item(265, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(265, 3).Value = "Bíobío"
$ws.Cells.Item(265, 4).Value = 45244
$ws.Cells.Item(265, 5).Value = 8
$ws.Cells.Item(265, 6).Value = "Fruta"
$ws.Cells.Item(265, 7).Value = 100108
$ws.Cells.Item(265, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(265, 9).Value = 100108005
$ws.Cells.Item(265, 10).Value = "Piña"
$ws.Cells.Item(265, 11).Value = "Caramelo"
$ws.Cells.Item(265, 12).Value = "Segunda"
$ws.Cells.Item(265, 13).Value = 200
$ws.Cells.Item(265, 14).Value = 20000
$ws.Cells.Item(265, 15).Value = 21000
$ws.Cells.Item(265, 16).Value = 20500
$ws.Cells.Item(265, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(265, 18).Value = "Ecuador"
$ws.Cells.Item(265, 19).Value = 1464
$ws.Cells.Item(265, 20).Value = 14
